$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new columns AD, AE, AF on row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of the existing header cell (AC1) to the new header cells
$ws.Range("AC1:AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill team record values for each data row (rows 2-43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "done"
